$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "Helliquary Raid"
$ws.Range("A6").Select()
